$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.257083535194397
$ws.Range("B1").Value = 1.429612159729004
$ws.Range("C1").Value = 3.683853626251221
$ws.Range("D1").Value = 3.520485401153564
$ws.Range("E1").Value = 1.008988261222839
